$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 2-21 down to rows 11-30 (old falling-event samples),
# iterating bottom-up so we never overwrite a row before reading it.
for ($r = 21; $r -ge 2; $r--) {
    $dest = $r + 9
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($dest, $c).Value = $ws.Cells.Item($r, $c).Value2
    }
}

# Insert the new samples captured on May 9th into rows 2-10.
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "falling"
$ws.Cells.Item(2,3).Value = -4.391921520233154
$ws.Cells.Item(2,4).Value = 5.980224609375
$ws.Cells.Item(2,5).Value = 2.74719500541687
$ws.Cells.Item(2,6).Value = 0.005824529005452867
$ws.Cells.Item(2,7).Value = -0.004876267489825652
$ws.Cells.Item(2,8).Value = 0.009258870057068522

$ws.Cells.Item(3,1).Value = 100
$ws.Cells.Item(3,2).Value = "falling"
$ws.Cells.Item(3,3).Value = -4.076038837432861
$ws.Cells.Item(3,4).Value = 6.086456298828125
$ws.Cells.Item(3,5).Value = 2.820330142974853
$ws.Cells.Item(3,6).Value = 0.01996676961696429
$ws.Cells.Item(3,7).Value = -0.05454103202494082
$ws.Cells.Item(3,8).Value = -0.001441926153939801

$ws.Cells.Item(4,1).Value = 200
$ws.Cells.Item(4,2).Value = "falling"
$ws.Cells.Item(4,3).Value = -3.987221717834473
$ws.Cells.Item(4,4).Value = 6.00103759765625
$ws.Cells.Item(4,5).Value = 2.755735635757446
$ws.Cells.Item(4,6).Value = -0.02569185483247746
$ws.Cells.Item(4,7).Value = -0.01999163068830969
$ws.Cells.Item(4,8).Value = -0.009539442396787617

$ws.Cells.Item(5,1).Value = 300
$ws.Cells.Item(5,2).Value = "falling"
$ws.Cells.Item(5,3).Value = -4.239626407623291
$ws.Cells.Item(5,4).Value = 5.927947998046875
$ws.Cells.Item(5,5).Value = 2.57724142074585
$ws.Cells.Item(5,6).Value = 0.03522419491999366
$ws.Cells.Item(5,7).Value = 0.005352173823603298
$ws.Cells.Item(5,8).Value = -0.01602810922317025

$ws.Cells.Item(6,1).Value = 400
$ws.Cells.Item(6,2).Value = "falling"
$ws.Cells.Item(6,3).Value = -4.282441139221191
$ws.Cells.Item(6,4).Value = 6.027595043182373
$ws.Cells.Item(6,5).Value = 2.64313268661499
$ws.Cells.Item(6,6).Value = -0.01577594932601883
$ws.Cells.Item(6,7).Value = 0.060614168860538
$ws.Cells.Item(6,8).Value = -0.009635333469960495

$ws.Cells.Item(7,1).Value = 500
$ws.Cells.Item(7,2).Value = "falling"
$ws.Cells.Item(7,3).Value = -4.193782329559326
$ws.Cells.Item(7,4).Value = 5.953823566436768
$ws.Cells.Item(7,5).Value = 2.714946508407593
$ws.Cells.Item(7,6).Value = -0.03060719080615873
$ws.Cells.Item(7,7).Value = 0.2389583984433218
$ws.Cells.Item(7,8).Value = -0.11525819691028

$ws.Cells.Item(8,1).Value = 600
$ws.Cells.Item(8,2).Value = "falling"
$ws.Cells.Item(8,3).Value = -4.003739356994629
$ws.Cells.Item(8,4).Value = 5.964433670043945
$ws.Cells.Item(8,5).Value = 2.773677349090576
$ws.Cells.Item(8,6).Value = -0.02314895105569856
$ws.Cells.Item(8,7).Value = 0.2151985930842024
$ws.Cells.Item(8,8).Value = -0.08631312587233481

$ws.Cells.Item(9,1).Value = 700
$ws.Cells.Item(9,2).Value = "falling"
$ws.Cells.Item(9,3).Value = -4.139037609100342
$ws.Cells.Item(9,4).Value = 5.919798851013184
$ws.Cells.Item(9,5).Value = 3.334548950195312
$ws.Cells.Item(9,6).Value = -0.06812567826966884
$ws.Cells.Item(9,7).Value = 0.2275579571723937
$ws.Cells.Item(9,8).Value = -0.07889750547880346

$ws.Cells.Item(10,1).Value = 800
$ws.Cells.Item(10,2).Value = "falling"
$ws.Cells.Item(10,3).Value = -4.075920104980469
$ws.Cells.Item(10,4).Value = 5.56472110748291
$ws.Cells.Item(10,5).Value = 3.763194084167481
$ws.Cells.Item(10,6).Value = -0.1647382801355317
$ws.Cells.Item(10,7).Value = 0.1269321128032929
$ws.Cells.Item(10,8).Value = -0.3623354202786154

# Append one more new sample as row 31.
$ws.Cells.Item(31,1).Value = 2900
$ws.Cells.Item(31,2).Value = "falling"
$ws.Cells.Item(31,3).Value = -0.2319450378417968
$ws.Cells.Item(31,4).Value = 3.570143699645996
$ws.Cells.Item(31,5).Value = 6.152892589569092
$ws.Cells.Item(31,6).Value = 0.02472228522217548
$ws.Cells.Item(31,7).Value = -0.03295831200341841
$ws.Cells.Item(31,8).Value = -0.03895686653464332
